# Seat Assignments sheet: move the row currently at A23:E23 down to the
# end of the contiguous block (new row 83), adding a Notes entry in H83.
#
# Net effect (matches the target diff):
#   - old row 23 content is removed from row 23; rows 24-83 shift up to 23-82
#   - a new row is opened at 83 and filled with the old row-23 values
#   - H83 gets a new "Notes" value
#   - row 84 and everything below is left completely untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# Capture the values currently stored in row 23 before we move anything.
$idVal          = $ws.Range("A23").Value2
$recordDayIdVal = $ws.Range("B23").Value2
$contestantVal  = $ws.Range("C23").Value2
$blockVal       = $ws.Range("D23").Value2
$seatVal        = $ws.Range("E23").Value2

# Delete row 23 entirely; rows below shift up by one.
$ws.Rows("23:23").Delete()

# Insert a fresh blank row at the position the moved row should end up
# (row 83); rows from 83 downward shift back down by one, restoring the
# original row numbering for row 84 onward.
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the captured values.
$ws.Range("A83").Value = $idVal
$ws.Range("B83").Value = $recordDayIdVal
$ws.Range("C83").Value = $contestantVal
$ws.Range("D83").Value = $blockVal
$ws.Range("E83").Value = $seatVal

# Add the new Notes value in column H for row 83.
$ws.Range("H83").Value = "idhsaiufahsdiufhasidufhuaidfhaisduhfaisuhfiuasdhfiuasdfiuasdhfiuadhsfiuadshfiuahsdifuahsdufuaishdfiudhf"
